# Revised abstract, methods, personnel
# Adds a new Personnel row (Harvey Walsh, technician, Northeast Fisheries
# Science Center / NOAA EcoMon) to the "Personnel" sheet, matching the
# existing style/format used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$wsPersonnel = $wb.Worksheets.Item("Personnel")
$wsKeywords  = $wb.Worksheets.Item("Keywords")

# New row 8 on the Personnel sheet:
#   A8 = Harvey            (given name)
#   C8 = Walsh             (surname)
#   D8 = Northeast Fisheries Science Center (organization)
#   G8 = technician        (role)
#   H8 = Northeast U.S. Shelf LTER (project title)
#   I8 = NOAA              (funding agency)
#   J8 = EcoMon            (funding number)
# B8, E8, F8 are intentionally left blank (no middle initial, email, or userId).

# Match formatting: A8/C8 use the same style as the name cells on the
# Keywords sheet (s="5"); D8/G8/H8/I8/J8 use the same style as the other
# data rows on the Personnel sheet (s="4", e.g. row 2).

$wsKeywords.Range("A1").Copy() | Out-Null
$wsPersonnel.Range("A8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("A8").Value = "Harvey"

$wsKeywords.Range("A1").Copy() | Out-Null
$wsPersonnel.Range("C8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("C8").Value = "Walsh"

$wsPersonnel.Range("D2").Copy() | Out-Null
$wsPersonnel.Range("D8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("D8").Value = "Northeast Fisheries Science Center"

$wsPersonnel.Range("G7").Copy() | Out-Null
$wsPersonnel.Range("G8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("G8").Value = "technician"

$wsPersonnel.Range("H2").Copy() | Out-Null
$wsPersonnel.Range("H8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("H8").Value = "Northeast U.S. Shelf LTER"

$wsPersonnel.Range("I3").Copy() | Out-Null
$wsPersonnel.Range("I8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("I8").Value = "NOAA"

$wsPersonnel.Range("J3").Copy() | Out-Null
$wsPersonnel.Range("J8").PasteSpecial(-4122) | Out-Null
$wsPersonnel.Range("J8").Value = "EcoMon"

$excel.CutCopyMode = 0
